# Update the "Riders" (C) and "Average" (D) columns on the Ridership sheet
# with the new Madigan bike hours figures.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Ridership")

$ws.Range("C2").Value = 286
$ws.Range("D2").Value = 235

$ws.Range("C3").Value = 241
$ws.Range("D3").Value = 211.64

$ws.Range("C4").Value = 200
$ws.Range("D4").Value = 195.29

$ws.Range("C5").Value = 196
$ws.Range("D5").Value = 222.46

$ws.Range("C6").Value = 257
$ws.Range("D6").Value = 238.64

$ws.Range("C7").Value = 107
$ws.Range("D7").Value = 120.2

$ws.Range("C8").Value = 43
$ws.Range("D8").Value = 101.08
